$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Worksheet rename & defined name follow automatically in Excel, but
#     make sure the autoFilter-backed defined name still resolves after
#     the rename (Excel keeps it in sync automatically). ---
$ws.Name = "Technology context"

# --- Row-level edits -------------------------------------------------
# Timeline gains two new, earlier entries (Geogebra 5.0 / BrailleR) that
# sort in before the existing MathJax-core-accessibility / Office-Math-
# Speech rows, plus a brand new entry at the bottom (Desmos).
#
# Shared-string append order matters for a byte-faithful rebuild, so the
# new "Desmos accessibility" row is filled in first (while it is still
# the last row of the table), and only afterwards do we insert the two
# rows above row 15 for Geogebra/BrailleR.

# 1) Append the Desmos accessibility row right after the current last
#    data row (16), reusing its formatting.
$ws.Range("A16:E16").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)
$ws.Range("C16").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("A17").Value = 43345
$ws.Range("C17").Value = "Desmos accessibility"
$ws.Range("D17").Value = "Student capability enhanced"
$ws.Range("E17").Value = "https://blog.desmos.com/articles/friday-fave-for-november-2/"

# 2) Insert two fresh rows above row 15 for the Geogebra / BrailleR
#    entries (this pushes the old row 15/16 -- and the Desmos row we
#    just added -- down by two).
$ws.Rows("15:16").Insert()

$ws.Range("A6:E6").Copy()
$ws.Range("A15:E16").PasteSpecial(-4122)
$ws.Range("B15:B16").Clear()

$ws.Range("A15").Value = 41885
$ws.Range("E15").Value = "https://wiki.geogebra.org/en/Accessibility"
$ws.Range("C15").Value = "Geogebra 5.0"
$ws.Range("D15").Value = "Student capability enhanced"

$ws.Range("A16").Value = 42094
$ws.Range("C16").Value = "BrailleR"
$ws.Range("D16").Value = "Student capability enhanced"
$ws.Range("E16").Value = "https://cran.r-project.org/web/packages/BrailleR/index.html"

# --- Selection / view tweaks -----------------------------------------
$ws.Range("B17").Select()
$excel.ActiveWindow.ScrollWorkbookTabs(1)
